$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  44"
$ws.Range("C9").Value = "Report Covering the Week  10/28/2024  Through  11/3/2024"

# --- Cells converting from text ("n/a"-style) to numeric ---
$ws.Range("C15").NumberFormat = '#,##0'
$ws.Range("C15").Value = 1
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("D15").Value = 1
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E15").Value = 0
$ws.Range("C22").NumberFormat = '#,##0'
$ws.Range("C22").Value = 1
$ws.Range("C23").NumberFormat = '#,##0'
$ws.Range("C23").Value = 1
$ws.Range("F23").NumberFormat = '#,##0'
$ws.Range("F23").Value = 1
$ws.Range("C27").NumberFormat = '#,##0'
$ws.Range("C27").Value = 1
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("D27").Value = 1
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E27").Value = 0
$ws.Range("C31").NumberFormat = '#,##0'
$ws.Range("C31").Value = 2
$ws.Range("D31").NumberFormat = '#,##0'
$ws.Range("D31").Value = 1
$ws.Range("E31").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E31").Value = 100
$ws.Range("F31").NumberFormat = '#,##0'
$ws.Range("F31").Value = 2

# --- Cells converting from numeric to text ("0" / "***.*") ---
# Reference cell that already carries the plain "text" style (s=13) we want to reuse
$textStyleSrc = $ws.Range("A14")
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0"
$textStyleSrc.Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "***.*"
$textStyleSrc.Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "0"
$textStyleSrc.Copy()
$ws.Range("G33").PasteSpecial(-4122)
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = "***.*"
$textStyleSrc.Copy()
$ws.Range("H33").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("M14").Value = -42.857142857142
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 29
$ws.Range("J15").Value = 33
$ws.Range("K15").Value = -12.121212121212
$ws.Range("L15").Value = -17.142857142857
$ws.Range("M15").Value = 107.142857142857
$ws.Range("N15").Value = 16
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = -80
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 27
$ws.Range("H16").Value = -51.851851851851
$ws.Range("I16").Value = 257
$ws.Range("J16").Value = 273
$ws.Range("K16").Value = -5.860805860805
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 48.554913294797
$ws.Range("N16").Value = -74.579624134520
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = -60
$ws.Range("F17").Value = 21
$ws.Range("G17").Value = 35
$ws.Range("H17").Value = -40
$ws.Range("I17").Value = 428
$ws.Range("J17").Value = 380
$ws.Range("K17").Value = 12.631578947368
$ws.Range("L17").Value = 50.175438596491
$ws.Range("M17").Value = 127.659574468085
$ws.Range("N17").Value = 40.327868852459
$ws.Range("C18").Value = 11
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 120
$ws.Range("F18").Value = 31
$ws.Range("G18").Value = 44
$ws.Range("H18").Value = -29.545454545454
$ws.Range("I18").Value = 345
$ws.Range("J18").Value = 468
$ws.Range("K18").Value = -26.282051282051
$ws.Range("L18").Value = -27.061310782241
$ws.Range("M18").Value = -16.058394160583
$ws.Range("N18").Value = -82.979773063640
$ws.Range("C19").Value = 23
$ws.Range("D19").Value = 21
$ws.Range("E19").Value = 9.523809523809
$ws.Range("F19").Value = 87
$ws.Range("G19").Value = 103
$ws.Range("H19").Value = -15.533980582524
$ws.Range("I19").Value = 1105
$ws.Range("J19").Value = 1160
$ws.Range("K19").Value = -4.741379310344
$ws.Range("L19").Value = -13.197172034564
$ws.Range("M19").Value = 95.229681978798
$ws.Range("N19").Value = -2.728873239436
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 32
$ws.Range("G20").Value = 56
$ws.Range("H20").Value = -42.857142857142
$ws.Range("I20").Value = 462
$ws.Range("J20").Value = 441
$ws.Range("K20").Value = 4.761904761904
$ws.Range("L20").Value = 102.631578947368
$ws.Range("M20").Value = 89.344262295082
$ws.Range("N20").Value = -86.666666666666
$ws.Range("C21").Value = 46
$ws.Range("D21").Value = 57
$ws.Range("E21").Value = -19.298245614035
$ws.Range("F21").Value = 186
$ws.Range("G21").Value = 266
$ws.Range("H21").Value = -30.075187969924
$ws.Range("I21").Value = 2630
$ws.Range("J21").Value = 2758
$ws.Range("K21").Value = -4.641044234952
$ws.Range("L21").Value = 2.935420743639
$ws.Range("M21").Value = 64.067373674360
$ws.Range("N21").Value = -67.038475999498
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -75
$ws.Range("I22").Value = 25
$ws.Range("J22").Value = 31
$ws.Range("K22").Value = -19.354838709677
$ws.Range("L22").Value = -30.555555555555
$ws.Range("M22").Value = 525
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = -66.666666666666
$ws.Range("I23").Value = 20
$ws.Range("J23").Value = 15
$ws.Range("K23").Value = 33.333333333333
$ws.Range("L23").Value = 53.846153846153
$ws.Range("M23").Value = 100
$ws.Range("C24").Value = 49
$ws.Range("D24").Value = 40
$ws.Range("E24").Value = 22.5
$ws.Range("F24").Value = 229
$ws.Range("G24").Value = 207
$ws.Range("H24").Value = 10.628019323671
$ws.Range("I24").Value = 2329
$ws.Range("J24").Value = 2385
$ws.Range("K24").Value = -2.348008385744
$ws.Range("L24").Value = -0.427533133817
$ws.Range("M24").Value = 73.031203566121
$ws.Range("C25").Value = 36
$ws.Range("D25").Value = 24
$ws.Range("E25").Value = 50
$ws.Range("F25").Value = 144
$ws.Range("G25").Value = 114
$ws.Range("H25").Value = 26.315789473684
$ws.Range("I25").Value = 1538
$ws.Range("J25").Value = 1309
$ws.Range("K25").Value = 17.494270435446
$ws.Range("L25").Value = 23.434991974317
$ws.Range("C26").Value = 21
$ws.Range("D26").Value = 22
$ws.Range("E26").Value = -4.545454545454
$ws.Range("F26").Value = 83
$ws.Range("G26").Value = 81
$ws.Range("H26").Value = 2.469135802469
$ws.Range("I26").Value = 828
$ws.Range("J26").Value = 762
$ws.Range("K26").Value = 8.661417322834
$ws.Range("L26").Value = 38.926174496644
$ws.Range("M26").Value = 35.960591133004
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 37
$ws.Range("J27").Value = 49
$ws.Range("K27").Value = -24.489795918367
$ws.Range("L27").Value = -19.565217391304
$ws.Range("C28").Value = 4
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 300
$ws.Range("F28").Value = 10
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = 42.857142857142
$ws.Range("I28").Value = 83
$ws.Range("J28").Value = 88
$ws.Range("K28").Value = -5.681818181818
$ws.Range("L28").Value = 22.058823529411
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 100
$ws.Range("I31").Value = 8
$ws.Range("J31").Value = 8
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = -38.461538461538
